$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "31.027.86"
$ws.Range("E2").Value = "  +1.46%  "
$ws.Range("D3").Value = "1.962.06"
$ws.Range("E3").Value = "  +2.26%  "
$ws.Range("D4").Value = "'1.0000"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'247.75"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "'0.4825"
$ws.Range("E7").Value = "  -0.71%  "
$ws.Range("D8").Value = "'0.2935"
$ws.Range("E8").Value = "  +1.09%  "
$ws.Range("D9").Value = "'0.06789"
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("D10").Value = "'19.14"
$ws.Range("E10").Value = "  -1.14%  "
$ws.Range("D11").Value = "'106.19"
$ws.Range("E11").Value = "  -4.91%  "
$ws.Range("D12").Value = "1.961.65"
$ws.Range("E12").Value = "  +2.10%  "
$ws.Range("D13").Value = "'0.07766"
$ws.Range("E13").Value = "  +2.34%  "
$ws.Range("D14").Value = "'5.447"
$ws.Range("E14").Value = "  +1.18%  "
$ws.Range("D15").Value = "'0.7023"
$ws.Range("E15").Value = "  +4.27%  "
$ws.Range("D16").Value = "'286.80"
$ws.Range("E16").Value = "  -2.95%  "
$ws.Range("D17").Value = "31.034.65"
$ws.Range("E17").Value = "  +1.48%  "
$ws.Range("D18").Value = "'13.20"
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("D19").Value = "'0.000007718"
$ws.Range("E19").Value = "  +2.01%  "
$ws.Range("D20").Value = "2.209.08"
$ws.Range("E20").Value = "  +1.67%  "
$ws.Range("D21").Value = "'5.588"
$ws.Range("E21").Value = "  +1.16%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "'6.585"
$ws.Range("E24").Value = "  +2.43%  "
$ws.Range("D25").Value = "'9.983"
$ws.Range("E25").Value = "  +5.17%  "
$ws.Range("D26").Value = "'169.94"
$ws.Range("E26").Value = "  +3.02%  "
$ws.Range("D27").Value = "'19.98"
$ws.Range("E27").Value = "  -2.13%  "
$ws.Range("D28").Value = "'2.185"
$ws.Range("E28").Value = "  +3.68%  "
$ws.Range("D29").Value = "'0.1061"
$ws.Range("E29").Value = "  -1.39%  "
$ws.Range("D30").Value = "'1.444"
$ws.Range("E30").Value = "  +0.64%  "
$ws.Range("D31").Value = "'4.800"
$ws.Range("E31").Value = "  +16.97%  "
$ws.Range("D32").Value = "'4.491"
$ws.Range("E32").Value = "  +8.40%  "
$ws.Range("D33").Value = "'0.05080"
$ws.Range("E33").Value = "  +1.27%  "
$ws.Range("D34").Value = "'0.7732"
$ws.Range("E34").Value = "  +4.28%  "
$ws.Range("D35").Value = "'1.171"
$ws.Range("E35").Value = "  +2.72%  "
$ws.Range("D36").Value = "'2.733"
$ws.Range("E36").Value = "  +1.09%  "
$ws.Range("D37").Value = "'0.02033"
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("D38").Value = "'2.720"
$ws.Range("E38").Value = "  +1.06%  "
$ws.Range("D39").Value = "'6.498"
$ws.Range("E39").Value = "  +10.41%  "
$ws.Range("D40").Value = "'2.118"
$ws.Range("E40").Value = "  +4.77%  "
$ws.Range("D41").Value = "'0.8911"
$ws.Range("E41").Value = "  +2.83%  "
$ws.Range("D42").Value = "'109.77"
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").Value = "'0.4446"
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("D44").Value = "'71.68"
$ws.Range("E44").Value = "  +2.48%  "
$ws.Range("D45").Value = "'1.001"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "'7.508"
$ws.Range("E46").Value = "  +3.33%  "
$ws.Range("D47").Value = "'9.462"
$ws.Range("E47").Value = "  +1.65%  "
$ws.Range("D48").Value = "'0.1270"
$ws.Range("E48").Value = "  +3.27%  "
$ws.Range("D49").Value = "'950.85"
$ws.Range("E49").Value = "  +11.26%  "
$ws.Range("D50").Value = "'35.93"
$ws.Range("E50").Value = "  +2.82%  "
$ws.Range("D51").Value = "'46.92"
$ws.Range("E51").Value = "  -3.01%  "
